# Update column F ("dSF") values for specific rows, per the repull/push of
# data and recalculated mean described in the commit message.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = -4
    8  = 0
    11 = -1
    26 = -1
    28 = 1
    33 = -5
    34 = 1
    36 = 0
    41 = -3
    45 = -9
    50 = -5
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
